$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 need to be added with the same look & feel (styles) as the
# existing data rows above them (rows 2-4). Row 4 already carries the
# correct per-column styles for columns A:J, so copy its formatting down
# into rows 5 and 6 (column K is left untouched - row 6 already has its own
# blank, styled K cell, and row 5 does not use column K) before writing the
# new values into them.
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)
$ws.Range("A4:J4").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)

# Column A first (GEF names), then column C (exporter names), so that new
# shared-string entries are interned in the same order Excel produced them.
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("C6").Value = "Scone exporter"

# Row 5 - "Crumpet" facility
$ws.Range("B5").Value = 20001371
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6 - "Scone" facility
$ws.Range("B6").Value = 20001371
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"
